$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before the row currently containing "franzosa_ControlvsCD_Fp"
# (row 8), and populate it with the franzosa_ControlvsCD_ConvCD data.
$ws.Rows.Item(8).Insert()
$ws.Cells.Item(8, 1).Value = "franzosa_ControlvsCD_ConvCD"
$ws.Cells.Item(8, 2).Value = 0
$ws.Cells.Item(8, 3).Value = 0
$ws.Cells.Item(8, 4).Value = 0
$ws.Cells.Item(8, 5).Value = 0.4
$ws.Cells.Item(8, 6).Value = 1
$ws.Cells.Item(8, 7).Value = 0.6
$ws.Cells.Item(8, 8).Value = 0.6

# Insert a new row before the row currently containing "franzosa_ControlvsUC_Fp"
# (now row 14 after the previous insertion), and populate it with the
# franzosa_ControlvsUC_ConvUC data.
$ws.Rows.Item(14).Insert()
$ws.Cells.Item(14, 1).Value = "franzosa_ControlvsUC_ConvUC"
$ws.Cells.Item(14, 2).Value = 0
$ws.Cells.Item(14, 3).Value = 0
$ws.Cells.Item(14, 4).Value = 0
$ws.Cells.Item(14, 5).Value = 0.3
$ws.Cells.Item(14, 6).Value = 1
$ws.Cells.Item(14, 7).Value = 0.7
$ws.Cells.Item(14, 8).Value = 0.7
